# ---------------------------------------------------------------------------
# Adds a new "bias adjustment" worksheet at the end of the workbook that
# consolidates the "Trainee sample" (treat=1, units 1-10) and the
# "Optimized sample" (treat=0) blocks into a single 5-column table, and
# updates a couple of stale sheet-view selections left over from editing.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Create the new sheet at the very end of the tab strip -------------
# Copying an existing sheet (rather than Worksheets.Add()) means the new
# sheet inherits the workbook's real default row height (16, i.e. the
# 12pt Calibri "Normal" style) instead of the generic 15pt default that a
# brand new blank sheet would get.
$sheet2 = $wb.Worksheets.Item("Sheet2")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet2.Copy($null, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "bias adjustment"

# Wipe the copied data/formatting - we only wanted the sheet's blank defaults.
$ws.Cells.Clear()

# --- 2. Header row ----------------------------------------------------------
$ws.Range("A1").Value = "treat"
$ws.Range("B1").Value = "unit"
$ws.Range("C1").Value = "age"
$ws.Range("D1").Value = "gpa"
$ws.Range("E1").Value = "earnings"

# --- 3. Trainee sample rows (treat = 1), units 1-10 -------------------------
$trainee = @(
  @(1, 18, 1.2768669668126549, 9500),
  @(2, 29, 2.8019828173782786, 12250),
  @(3, 24, 3.9205642309579938, 11000),
  @(4, 27, 2.2926678508031757, 11750),
  @(5, 33, 2.4966702622946206, 13250),
  @(6, 22, 1.3387078036620916, 10500),
  @(7, 19, 1.662759442733091, 9750),
  @(8, 20, 2.598744835539557, 10000),
  @(9, 21, 1.9413317247167639, 10250),
  @(10, 30, 3.3740289080760242, 12500)
)
for ($i = 0; $i -lt $trainee.Count; $i++) {
  $r = 2 + $i
  $row = $trainee[$i]
  $ws.Range("A" + $r).Value = 1
  $ws.Range("B" + $r).Value = $row[0]
  $ws.Range("C" + $r).Value = $row[1]
  $ws.Range("D" + $r).Value = $row[2]
  $ws.Range("E" + $r).Value = $row[3]
}

# --- 4. Comparison / optimized sample rows (treat = 0) ----------------------
$comparison = @(
  @(13, 22, 1.662291385973246, 8950),
  @(5, 38, 1.6145408563338297, 12550),
  @(8, 33, 1.9747705998072842, 11425),
  @(2, 27, 1.77610720900872, 10075),
  @(8, 33, 1.9747705998072842, 11425),
  @(13, 22, 1.662291385973246, 8950),
  @(17, 19, 1.8595863891608966, 8275),
  @(1, 20, 1.8904439599777811, 8500),
  @(3, 21, 1.837511566218393, 8725),
  @(10, 30, 2.0152083863823114, 10750)
)
for ($i = 0; $i -lt $comparison.Count; $i++) {
  $r = 12 + $i
  $row = $comparison[$i]
  $ws.Range("A" + $r).Value = 0
  $ws.Range("B" + $r).Value = $row[0]
  $ws.Range("C" + $r).Value = $row[1]
  $ws.Range("D" + $r).Value = $row[2]
  $ws.Range("E" + $r).Value = $row[3]
}

# --- 5. Reuse the existing visual formatting (medium-border boxes + 2dp
#        gpa number format) from the tables these numbers were copied
#        from, instead of re-building ad-hoc borders. ----------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("A2:D12").Copy()
$ws.Range("B1:E11").PasteSpecial(-4122)

$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet3.Range("I3:L12").Copy()
$ws.Range("B12:E21").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 6. Sheet view bits: new sheet becomes the active tab, selection at F1 -
$ws.Range("F1").Select() | Out-Null
$ws.Activate() | Out-Null

# --- 7. Stale selection left on "Sheet3" (the sheet with the Optimized
#        sample block) gets updated to cover that block. -------------------
$sheet3.Range("I3:L12").Select() | Out-Null

$ws.Activate() | Out-Null
